# Auto-generated: updates cryptos list price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.958.08'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').Value = '2.335.43'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.536'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = '2.370.07'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.102'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('E11').Value = '  +2.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.355'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.45%  '
$ws.Range('D14').Value = '2.770.93'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = '57.954.23'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = '2.364.01'
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '340.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.42%  '
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('E21').Value = '  +1.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('E25').Value = '  +2.80%  '
$ws.Range('E26').Value = '  -1.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').Value = '  +6.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '175.02'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.92%  '
$ws.Range('E30').Value = '  +4.89%  '
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('E34').Value = '  +13.18%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.992'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.27'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('E38').Value = '  +3.31%  '
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('E40').Value = '  +2.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '150.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '284.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0506'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '18.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.56%  '
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.01%  '
$ws.Range('E51').Value = '  +6.62%  '
